$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71; this shifts existing rows 71..142 down to 72..143
# and copies formatting (e.g. the date style on column D) down from the row above.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new record's data.
$ws.Range("A71").Value = 1
$ws.Range("B71").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C71").Value = "Arica y Parinacota"
$ws.Range("D71").Value = 44484
$ws.Range("E71").Value = 15
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100108
$ws.Range("H71").Value = "Tropicales y subtropicales"
$ws.Range("I71").Value = 100108006
$ws.Range("J71").Value = "Plátano"
$ws.Range("K71").Value = "Sin especificar"
$ws.Range("L71").Value = "Pintón"
$ws.Range("M71").Value = 120
$ws.Range("N71").Value = 21000
$ws.Range("O71").Value = 22000
$ws.Range("P71").Value = 21500
$ws.Range("Q71").Value = "`$/caja 20 kilos"
$ws.Range("R71").Value = "Bolivia"
$ws.Range("S71").Value = 1075
$ws.Range("T71").Value = 20
